# Update theta_se (row 4) and lambda_se (row 6) standard-error values
# to reflect the new bootstrapping results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: theta_se
$ws.Range("B4").Value = "(0.05)"
$ws.Range("C4").Value = "(0.18)"
$ws.Range("D4").Value = "(0.27)"
$ws.Range("E4").Value = "(0.08)"
$ws.Range("F4").Value = "(0.62)"
$ws.Range("G4").Value = "(0.2)"

# Row 6: lambda_se
$ws.Range("B6").Value = "(0.04)"
$ws.Range("C6").Value = "(0.12)"
$ws.Range("D6").Value = "(0.23)"
$ws.Range("E6").Value = "(0.08)"
$ws.Range("F6").Value = "(0.23)"
$ws.Range("G6").Value = "(0.4)"
